# Apply cryptocurrency price/volume updates from the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted numbers (e.g. "3.061.81") in the source
# data, so force Text number format before assigning to avoid Excel silently
# re-typing these as numeric values (which would drop the thousands-dot style).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D25", "D28", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($cell in $priceCells) { $ws.Range($cell).NumberFormat = "@" }

$ws.Range("D2").Value = '63.942.55'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '3.064.67'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '559.02'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = '143.05'
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D8").Value = '3.063.15'
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("E9").Value = '  +3.63%  '
$ws.Range("D10").Value = '0.154'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").Value = '6.16'
$ws.Range("E11").Value = '  -3.98%  '
$ws.Range("D12").Value = '0.482'
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").Value = '0.0000233'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").Value = '35.41'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '3.566.93'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '63.994.76'
$ws.Range("E16").Value = '  -1.41%  '
$ws.Range("D17").Value = '3.060.80'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '6.80'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = '487.87'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").Value = '14.37'
$ws.Range("E21").Value = '  +3.67%  '
$ws.Range("D22").Value = '0.693'
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("E23").Value = '  +8.08%  '
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = '82.78'
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").Value = '26.57'
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").Value = '5.71'
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("D35").Value = '6.25'
$ws.Range("E35").Value = '  +1.53%  '
$ws.Range("D36").Value = '54.86'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").Value = '0.0413'
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").Value = '444.80'
$ws.Range("E38").Value = '  -6.11%  '
$ws.Range("D39").Value = '0.0815'
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("D40").Value = '3.039.80'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").Value = '2.79'
$ws.Range("E41").Value = '  -5.88%  '
$ws.Range("D42").Value = '8.35'
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("E43").Value = '  +1.74%  '
$ws.Range("D44").Value = '0.276'
$ws.Range("E44").Value = '  +6.59%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '27.97'
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.28'
$ws.Range("E46").Value = '  +6.30%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '0.114'
$ws.Range("E48").Value = '  +1.26%  '
$ws.Range("D49").Value = '0.0₃0519'
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("D50").Value = '117.93'
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("D51").Value = '2.14'
$ws.Range("E51").Value = '  +3.29%  '
